$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")
Write-Host $ws.Name
